$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-09 12:41:49"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
